$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix error in abbreviation row: replace "NA" placeholders with actual numeric values
$ws.Range("B11").Value = 1525
$ws.Range("C11").Value = 52811
$ws.Range("D11").Value = 1661
